$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: headers for default_next_levelup_formula block
$ws.Range("A8").Value = "default_next_levelup_formula"
$ws.Range("B8").Value = "level"
$ws.Range("C8").Value = "next levelup defeats"

# Row 9: values/formula
$ws.Range("B9").Value = 16
$ws.Range("C9").Formula = "=FLOOR(1+(B9*B9*0.25),1)"

# Row 10: headers for classic_next_levelup_formula block
$ws.Range("A10").Value = "classic_next_levelup_formula"
$ws.Range("B10").Value = "level"
$ws.Range("C10").Value = "next levelup defeats"

# Row 11: values/formula
$ws.Range("B11").Value = 16
$ws.Range("C11").Formula = "=FLOOR((2+B11)*0.7,1)"

# Column C width similar to column E (bestFit-like width)
$ws.Columns.Item(3).ColumnWidth = 8.65

# Selection matches final state in diff
$ws.Range("C11").Select()
